# Add two new weekly reporting periods ("7.-13. 6. 2021" and "14.-20. 6. 2021")
# to both the "data" and "pocetR" sheets, and bump the "aktualizace" date in the
# two footer/title cells from "1. 6. 2021" to "28. 6. 2021".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "data": new columns BF (period "7.-13. 6. 2021") and BG (period
# "14.-20. 6. 2021") appended right after the existing BE ("17.-23. 5. 2021").
# ---------------------------------------------------------------------------
$wsData = $wb.Worksheets.Item("data")

# Carry over the header formatting (bold, centered, bordered) from BE1.
$wsData.Range("BE1").Copy()
$wsData.Range("BF1:BG1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsData.Range("BF1").Value = "7.–13. 6. 2021"
$wsData.Range("BG1").Value = "14.–20. 6. 2021"

$dataNewCols = @{
    2 = @(0.75, 0.76)
    3 = @(0.08, 0.08)
    4 = @(0.09, 0.08)
    5 = @(0.08, 0.08)
    6 = @(0.74, 0.76)
    7 = @(0.08, 0.07000000000000001)
    8 = @(0.09, 0.08)
    9 = @(0.09, 0.09)
    10 = @(0.45, 0.47)
    11 = @(0.19, 0.17)
    12 = @(0.36, 0.34)
    13 = @(0, 0.02)
    14 = @(0.83, 0.84)
    15 = @(0.05, 0.08)
    16 = @(0.04, 0.04)
    17 = @(0.08, 0.04)
    18 = @(0.9, 0.86)
    19 = @(0.03, 0.05)
    20 = @(0.03, 0.03)
    21 = @(0.04, 0.06)
    22 = @(0.72, 0.76)
    23 = @(0.16, 0.14)
    24 = @(0.03, 0.015)
    25 = @(0.09, 0.08500000000000001)
    26 = @(0.8, 0.8100000000000001)
    27 = @(0.06, 0.07000000000000001)
    28 = @(0.07000000000000001, 0.06)
    29 = @(0.07000000000000001, 0.06)
    30 = @(0.79, 0.8)
    31 = @(0.06, 0.06)
    32 = @(0.06, 0.06)
    33 = @(0.09, 0.08)
    34 = @(0.66, 0.71)
    35 = @(0.13, 0.09)
    36 = @(0.14, 0.13)
    37 = @(0.07000000000000001, 0.07000000000000001)
    38 = @(0.62, 0.62)
    39 = @(0.16, 0.14)
    40 = @(0.16, 0.16)
    41 = @(0.06, 0.08)
    42 = @(0.71, 0.72)
    43 = @(0.11, 0.12)
    44 = @(0.13, 0.14)
    45 = @(0.05, 0.02)
    46 = @(0.79, 0.8100000000000001)
    47 = @(0.07000000000000001, 0.06)
    48 = @(0.07000000000000001, 0.07000000000000001)
    49 = @(0.07000000000000001, 0.06)
    50 = @(0.7, 0.6899999999999999)
    51 = @(0.07000000000000001, 0.09)
    52 = @(0.07000000000000001, 0.07000000000000001)
    53 = @(0.16, 0.15)
    54 = @(0.76, 0.8)
    55 = @(0.09, 0.07000000000000001)
    56 = @(0.05, 0.05)
    57 = @(0.1, 0.08)
    58 = @(0.72, 0.72)
    59 = @(0.07000000000000001, 0.1)
    60 = @(0.12, 0.11)
    61 = @(0.09, 0.07000000000000001)
    62 = @(0.76, 0.76)
    63 = @(0.08, 0.08)
    64 = @(0.09, 0.09)
    65 = @(0.07000000000000001, 0.07000000000000001)
    66 = @(0.8, 0.83)
    67 = @(0.05, 0.05)
    68 = @(0.03, 0.04)
    69 = @(0.12, 0.08)
    70 = @(0.75, 0.76)
    71 = @(0.09, 0.09)
    72 = @(0.08, 0.08)
    73 = @(0.08, 0.07000000000000001)
    74 = @(0.67, 0.67)
    75 = @(0.12, 0.12)
    76 = @(0.19, 0.17)
    77 = @(0.02, 0.04)
}

foreach ($row in ($dataNewCols.Keys | Sort-Object)) {
    $vals = $dataNewCols[$row]
    $wsData.Cells.Item($row, 58).Value = $vals[0]
    $wsData.Cells.Item($row, 59).Value = $vals[1]
}

# Footer cell A78 carries the report title / refresh date - bump it.
$oldTitle = $wsData.Cells.Item(78, 1).Text
$newTitle = $oldTitle.Replace("1. 6. 2021", "28. 6. 2021")
$wsData.Cells.Item(78, 1).Value = $newTitle

# ---------------------------------------------------------------------------
# Sheet "pocetR": new columns BE (period "7.-13. 6. 2021") and BF (period
# "14.-20. 6. 2021") appended right after the existing BD ("17.-23. 5. 2021").
# ---------------------------------------------------------------------------
$wsPocetR = $wb.Worksheets.Item("pocetR")

$wsPocetR.Range("BD1").Copy()
$wsPocetR.Range("BE1:BF1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsPocetR.Range("BE1").Value = "7.–13. 6. 2021"
$wsPocetR.Range("BF1").Value = "14.–20. 6. 2021"

$pocetRNewCols = @{
    2 = @(1059, 1059)
    3 = @(302, 302)
    4 = @(99, 99)
    5 = @(268, 268)
    6 = @(144, 144)
    7 = @(95, 95)
    8 = @(501, 501)
    9 = @(267, 267)
    10 = @(135, 135)
    11 = @(156, 156)
    12 = @(256, 256)
    13 = @(600, 600)
    14 = @(203, 203)
    15 = @(269, 269)
    16 = @(214, 214)
    17 = @(576, 576)
    18 = @(385, 385)
    19 = @(422, 422)
    20 = @(252, 252)
}

foreach ($row in ($pocetRNewCols.Keys | Sort-Object)) {
    $vals = $pocetRNewCols[$row]
    $wsPocetR.Cells.Item($row, 57).Value = $vals[0]
    $wsPocetR.Cells.Item($row, 58).Value = $vals[1]
}

# Trailing blank-row formatting (row 21) is extended with two more blank,
# text-formatted cells, matching the rest of the row.
$wsPocetR.Range("BD21").Copy()
$wsPocetR.Range("BE21:BF21").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Footer cell A21 carries the report title / refresh date - bump it.
$oldTitle2 = $wsPocetR.Cells.Item(21, 1).Text
$newTitle2 = $oldTitle2.Replace("1. 6. 2021", "28. 6. 2021")
$wsPocetR.Cells.Item(21, 1).Value = $newTitle2

Write-Output "Done updating data/pocetR with the 2021-06-28 refresh (2 new weekly columns + title dates)."
